$wb = $excel.ActiveWorkbook
$wsTypo = $wb.Worksheets.Item('Typography')
$wsTrans = $wb.Worksheets.Item('Translation')

# --- Typography sheet updates (row 4 and row 6 wildcard settings) ---
$wsTypo.Range('G4').Value = '.",°βα'
$wsTypo.Range('I4').Value = '0-9,'
$wsTypo.Range('I6').Value = '0x0020-0x007F,0x00C0-0x00FF,0x0018,0x000D'

# --- Translation sheet: update existing rows (encoder/radius labels, feed rate units) ---
$wsTrans.Range('F75').Value = '   rez [mm/min]:'
$wsTrans.Range('G75').Value = 'f. rate [mm/min]:'
$wsTrans.Range('F77').Value = 'Shrani 
1'
$wsTrans.Range('G77').Value = 'Save 
1'
$wsTrans.Range('F78').Value = 'Shrani
2'
$wsTrans.Range('G78').Value = 'Save 
2'
$wsTrans.Range('F79').Value = 'Nalozi
1'
$wsTrans.Range('G79').Value = 'Load 
1'
$wsTrans.Range('F80').Value = 'Nalozi
2'
$wsTrans.Range('G80').Value = 'Load 
2'
$wsTrans.Range('F81').Value = 'Shrani
3'
$wsTrans.Range('G81').Value = 'Save 
3'
$wsTrans.Range('F82').Value = 'Nalozi
3'
$wsTrans.Range('G82').Value = 'Load 
3'
$wsTrans.Range('F85').Value = 'hitrost: <feedrate>mm/min'
$wsTrans.Range('G85').Value = 'f. rate: <feedrate>mm/min'

# --- Translation sheet: new rows 98-107 (GRBL / encoder / laser console strings) ---
$wsTrans.Range('B98').Value = 'SingleUseId106'
$wsTrans.Range('C98').Value = 'Large'
$wsTrans.Range('D98').Value = 'Left'
$wsTrans.Range('E98').Value = 'LTR'
$wsTrans.Range('F98').Value = 'mm'
$wsTrans.Range('G98').Value = 'mm'

$wsTrans.Range('B99').Value = 'SingleUseId107'
$wsTrans.Range('C99').Value = 'Default'
$wsTrans.Range('D99').Value = 'Right'
$wsTrans.Range('E99').Value = 'LTR'
$wsTrans.Range('F99').Value = 'Obseg'
$wsTrans.Range('G99').Value = 'Circumference'

$wsTrans.Range('B100').Value = 'SingleUseId108'
$wsTrans.Range('C100').Value = 'Large'
$wsTrans.Range('D100').Value = 'Left'
$wsTrans.Range('E100').Value = 'LTR'
$wsTrans.Range('F100').Value = '.'
$wsTrans.Range('G100').Value = '.'

$wsTrans.Range('B101').Value = 'SingleUseId109'
$wsTrans.Range('C101').Value = 'Default'
$wsTrans.Range('D101').Value = 'Right'
$wsTrans.Range('E101').Value = 'LTR'
$wsTrans.Range('F101').Value = 'Enkoder
p.n.o.'
$wsTrans.Range('G101').Value = 'Encoder
p.p.r'

$wsTrans.Range('B102').Value = 'SingleUseId110'
$wsTrans.Range('C102').Value = 'Default'
$wsTrans.Range('D102').Value = 'Center'
$wsTrans.Range('E102').Value = 'LTR'
$wsTrans.Range('F102').Value = 'Uporabi
radij'
$wsTrans.Range('G102').Value = 'Use 
radius'

$wsTrans.Range('B103').Value = 'STATUSMSG_ORIGIN_UPDATED'
$wsTrans.Range('C103').Value = 'Default'
$wsTrans.Range('D103').Value = 'Left'
$wsTrans.Range('E103').Value = 'LTR'
$wsTrans.Range('F103').Value = 'Izhodisce in domaca pozicija posodobljeni.'
$wsTrans.Range('G103').Value = 'Origin and home position updated.'

$wsTrans.Range('B104').Value = 'SingleUseId111'
$wsTrans.Range('C104').Value = 'Large'
$wsTrans.Range('D104').Value = 'Left'
$wsTrans.Range('E104').Value = 'LTR'
$wsTrans.Range('F104').Value = 'proiz.: <feedrate>mm/min'
$wsTrans.Range('G104').Value = ' prod.: <velocity>mm/min'

$wsTrans.Range('B105').Value = 'SingleUseId112'
$wsTrans.Range('C105').Value = 'Default'
$wsTrans.Range('D105').Value = 'Center'
$wsTrans.Range('E105').Value = 'LTR'
$wsTrans.Range('F105').Value = 'GRBL kontroler povezan. 
Nadaljujem z iskanjem (0,0)?'
$wsTrans.Range('G105').Value = 'GRBL controller connected. 
Proceed with homing?'

$wsTrans.Range('B106').Value = 'SingleUseId114'
$wsTrans.Range('C106').Value = 'Default'
$wsTrans.Range('D106').Value = 'Center'
$wsTrans.Range('E106').Value = 'LTR'
$wsTrans.Range('F106').Value = 'α rez'
$wsTrans.Range('G106').Value = 'α cut'

$wsTrans.Range('B107').Value = 'SingleUseId115'
$wsTrans.Range('C107').Value = 'Default'
$wsTrans.Range('D107').Value = 'Center'
$wsTrans.Range('E107').Value = 'LTR'
$wsTrans.Range('F107').Value = 'Laser 
Konzola'
$wsTrans.Range('G107').Value = 'Laser
Console'

